$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = -0.0736143064681548
$ws.Range("D2").Value = 0.2023283625086515
$ws.Range("E2").Value = 0.005361901149070607
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0.005910374655943606
$ws.Range("M2").Value = 0.2017405004068997
$ws.Range("N2").Value = -0.003175673222564392
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = -0.0836551046379089
$ws.Range("V2").Value = 0.01416007321150566
$ws.Range("W2").Value = -0.03298341659304817
$ws.Range("AB2").Value = 0
$ws.Range("AC2").Value = -0.05416417955287071
$ws.Range("AE2").Value = -0.01788162495550331
$ws.Range("AF2").Value = 0.0002069792777307436
$ws.Range("AK2").Value = 0
$ws.Range("AL2").Value = -0.03184932082569965
$ws.Range("AN2").Value = 0.02827770634814052
$ws.Range("AO2").Value = 0.06944358562979185
$ws.Range("AT2").Value = 0
$ws.Range("AU2").Value = -0.1497027310705481
$ws.Range("AW2").Value = 0.07080831603100772
$ws.Range("AX2").Value = -0.001817919973327277
$ws.Range("BC2").Value = 0
$ws.Range("BD2").Value = -0.01418805710578807
$ws.Range("BF2").Value = 0.08649163433815991
$ws.Range("BG2").Value = 0.03283123518905573
$ws.Range("BL2").Value = 0
$ws.Range("BM2").Value = 0.03190481457958391
$ws.Range("BO2").Value = -0.04257541708426302
$ws.Range("BP2").Value = -0.08905310676590357
$ws.Range("BU2").Value = 0
$ws.Range("BV2").Value = -0.04640710802875297
$ws.Range("BX2").Value = 0.01094127294829059
$ws.Range("BY2").Value = -0.02004983166574863
$ws.Range("CD2").Value = 0
$ws.Range("CE2").Value = 0.03274941356648393
$ws.Range("CG2").Value = -0.03099117283538349
$ws.Range("CH2").Value = 0.01593166186343111
$ws.Range("CM2").Value = 0
$ws.Range("CN2").Value = -0.01061381960660221
$ws.Range("CP2").Value = 0.02133540246658532
$ws.Range("CQ2").Value = 0.03716094318380431
$ws.Range("CV2").Value = 0
$ws.Range("CW2").Value = 0.04597314766486385
$ws.Range("CY2").Value = -0.03342648399499332
$ws.Range("CZ2").Value = 0.01017704690408558
$ws.Range("DE2").Value = 0
$ws.Range("DF2").Value = 0.02874569132567836
$ws.Range("DH2").Value = 0.02900526664094873
$ws.Range("DI2").Value = 0.03380050877759293
$ws.Range("DN2").Value = 0
$ws.Range("DO2").Value = -0.01950247745448723
$ws.Range("DQ2").Value = 0.03637034262361485
$ws.Range("DR2").Value = -0.01945341551444906
$ws.Range("DW2").Value = 0
$ws.Range("DX2").Value = -0.05668836815106189
$ws.Range("DZ2").Value = -0.008248668484950638
$ws.Range("EA2").Value = -0.02434073422596091
$ws.Range("EF2").Value = 0
$ws.Range("EG2").Value = 0.04085971144248264
$ws.Range("EI2").Value = 0.06686095049629477
$ws.Range("EJ2").Value = -0.02407360759003618
$ws.Range("EO2").Value = 0
$ws.Range("EP2").Value = 0.04575396385905522
$ws.Range("ER2").Value = -0.0349068518066118
$ws.Range("ES2").Value = 0.03654119765287879
$ws.Range("EX2").Value = 0
$ws.Range("EY2").Value = 0.04340351386436194
$ws.Range("FA2").Value = -0.02669773779825179
$ws.Range("FB2").Value = 0.01710822948871973
$ws.Range("FG2").Value = 0
$ws.Range("FH2").Value = 0.001691558334483329
$ws.Range("FJ2").Value = -0.006288705109680439
$ws.Range("FK2").Value = -0.008476023910300627
$ws.Range("FP2").Value = 0
$ws.Range("FQ2").Value = -0.01406165917757605
$ws.Range("FS2").Value = -0.0184191401413391
$ws.Range("FT2").Value = 0.005010644615690384
$ws.Range("FY2").Value = 0
$ws.Range("FZ2").Value = -0.03040407898268115
$ws.Range("GB2").Value = 0.03399166782388464
